$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.726732730865479
$ws.Range("B1").Value = 3.525825500488281
$ws.Range("C1").Value = 6.394368648529053
$ws.Range("D1").Value = 1.696390628814697
$ws.Range("E1").Value = 0.8720933794975281
